$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.9589784145355225
$ws.Range("E2").Value = 325.6107404475042
$ws.Range("F2").Value = 0.01522295381554292
$ws.Range("G2").Value = 0.01203426725805245
$ws.Range("H2").Value = 0.01046367017614035
$ws.Range("I2").Value = 0.009387097712417695
$ws.Range("J2").Value = 0.008743236199409411
$ws.Range("K2").Value = 0.007874141587312228
$ws.Range("L2").Value = 0.007685012470354559
$ws.Range("M2").Value = 0.007466545301897488
$ws.Range("N2").Value = 0.00725857159268521
$ws.Range("O2").Value = 0.007028698450637365
$ws.Range("P2").Value = 0.007001329766525526
$ws.Range("Q2").Value = 0.007000897497714743
$ws.Range("R2").Value = 0.006903842386797199
$ws.Range("S2").Value = 0.006823619127550031
$ws.Range("T2").Value = 0.006570010897062495
$ws.Range("U2").Value = 0.006550011282277756
$ws.Range("V2").Value = 0.006434373789056903
$ws.Range("W2").Value = 0.006434373789056903
$ws.Range("X2").Value = 0.006403016657405585
$ws.Range("Y2").Value = 0.006347187922953297
$ws.Range("C3").Value = 0.8359975814819336
$ws.Range("E3").Value = 323.4395295460508
$ws.Range("F3").Value = 0.01520188765163886
$ws.Range("G3").Value = 0.01162056397092288
$ws.Range("H3").Value = 0.009358328249736177
$ws.Range("I3").Value = 0.008506827521800124
$ws.Range("J3").Value = 0.008206739437993486
$ws.Range("K3").Value = 0.007529872942344929
$ws.Range("L3").Value = 0.007264489094723102
$ws.Range("M3").Value = 0.007077128908558042
$ws.Range("N3").Value = 0.006889356920258995
$ws.Range("O3").Value = 0.006830762990116426
$ws.Range("P3").Value = 0.006698157538724481
$ws.Range("Q3").Value = 0.006414141173685222
$ws.Range("R3").Value = 0.006414141173685222
$ws.Range("S3").Value = 0.006414141173685222
$ws.Range("T3").Value = 0.006414141173685222
$ws.Range("U3").Value = 0.006414141173685222
$ws.Range("V3").Value = 0.006381795563233501
$ws.Range("W3").Value = 0.006381795563233501
$ws.Range("X3").Value = 0.006366880323735583
$ws.Range("Y3").Value = 0.006304864123704693
$ws.Range("C4").Value = 0.8600001335144043
$ws.Range("E4").Value = 326.3546697341662
$ws.Range("F4").Value = 0.01508564402090302
$ws.Range("G4").Value = 0.01123366606423435
$ws.Range("H4").Value = 0.009438396495582363
$ws.Range("I4").Value = 0.009054009278655012
$ws.Range("J4").Value = 0.008764869378956772
$ws.Range("K4").Value = 0.008012932115079589
$ws.Range("L4").Value = 0.007632157328741571
$ws.Range("M4").Value = 0.007447872380018795
$ws.Range("N4").Value = 0.007096562686316356
$ws.Range("O4").Value = 0.00700875915653797
$ws.Range("P4").Value = 0.006830646056521721
$ws.Range("Q4").Value = 0.006763443741447136
$ws.Range("R4").Value = 0.006624143401553502
$ws.Range("S4").Value = 0.006557537439805536
$ws.Range("T4").Value = 0.006524634158498134
$ws.Range("U4").Value = 0.006408832913010036
$ws.Range("V4").Value = 0.006408832913010036
$ws.Range("W4").Value = 0.006408832913010036
$ws.Range("X4").Value = 0.006361689468502265
$ws.Range("Y4").Value = 0.006361689468502265
$ws.Range("C5").Value = 0.8229987621307373
$ws.Range("E5").Value = 329.3532944061844
$ws.Range("F5").Value = 0.01533738425161548
$ws.Range("G5").Value = 0.01242650990702754
$ws.Range("H5").Value = 0.01104618230891453
$ws.Range("I5").Value = 0.009895330762621354
$ws.Range("J5").Value = 0.009233906843296526
$ws.Range("K5").Value = 0.008263362854982421
$ws.Range("L5").Value = 0.008227962194479103
$ws.Range("M5").Value = 0.008037964958018161
$ws.Range("N5").Value = 0.00769237197822747
$ws.Range("O5").Value = 0.007145283092498087
$ws.Range("P5").Value = 0.006987939715046188
$ws.Range("Q5").Value = 0.006941621845703708
$ws.Range("R5").Value = 0.006716562396827017
$ws.Range("S5").Value = 0.006647970373464894
$ws.Range("T5").Value = 0.006647970373464894
$ws.Range("U5").Value = 0.006647970373464894
$ws.Range("V5").Value = 0.00652946908296456
$ws.Range("W5").Value = 0.006471545001469941
$ws.Range("X5").Value = 0.006420142191153692
$ws.Range("Y5").Value = 0.006420142191153692
$ws.Range("C6").Value = 0.7529628276824951
$ws.Range("E6").Value = 333.719393643536
$ws.Range("F6").Value = 0.01527701109838837
$ws.Range("G6").Value = 0.01199372904564168
$ws.Range("H6").Value = 0.01086043788878272
$ws.Range("I6").Value = 0.009446929706548641
$ws.Range("J6").Value = 0.008367819574422253
$ws.Range("K6").Value = 0.007936249144742304
$ws.Range("L6").Value = 0.007665486604341226
$ws.Range("M6").Value = 0.007400980125594338
$ws.Range("N6").Value = 0.007159706820272055
$ws.Range("O6").Value = 0.00700372480019904
$ws.Range("P6").Value = 0.006884033529429402
$ws.Range("Q6").Value = 0.006674887404769903
$ws.Range("R6").Value = 0.006674887404769903
$ws.Range("S6").Value = 0.006674887404769903
$ws.Range("T6").Value = 0.006641558733304459
$ws.Range("U6").Value = 0.00658818246601971
$ws.Range("V6").Value = 0.00657714672163343
$ws.Range("W6").Value = 0.006535129260011625
$ws.Range("X6").Value = 0.006532107874889715
$ws.Range("Y6").Value = 0.006505251338080623
$ws.Range("C7").Value = 0.8060019016265869
$ws.Range("E7").Value = 334.3752094244355
$ws.Range("F7").Value = 0.01527005842666242
$ws.Range("G7").Value = 0.01253699458197416
$ws.Range("H7").Value = 0.01026905507399499
$ws.Range("I7").Value = 0.009251222290491979
$ws.Range("J7").Value = 0.009093260832317089
$ws.Range("K7").Value = 0.008471349291259547
$ws.Range("L7").Value = 0.007721843413224705
$ws.Range("M7").Value = 0.007721843413224705
$ws.Range("N7").Value = 0.007503179756913949
$ws.Range("O7").Value = 0.007503179756913949
$ws.Range("P7").Value = 0.00722720924085626
$ws.Range("Q7").Value = 0.006975690658318976
$ws.Range("R7").Value = 0.006975690658318976
$ws.Range("S7").Value = 0.00692788882825091
$ws.Range("T7").Value = 0.006785971799456117
$ws.Range("U7").Value = 0.006738721039002508
$ws.Range("V7").Value = 0.006645004710487988
$ws.Range("W7").Value = 0.006580999843249617
$ws.Range("X7").Value = 0.006546625533174291
$ws.Range("Y7").Value = 0.006518035271431491
$ws.Range("C8").Value = 0.7199945449829102
$ws.Range("E8").Value = 322.8910549316734
$ws.Range("F8").Value = 0.01546014076994576
$ws.Range("G8").Value = 0.01253552032953864
$ws.Range("H8").Value = 0.01044941930156324
$ws.Range("I8").Value = 0.009703922213497024
$ws.Range("J8").Value = 0.009158192046421073
$ws.Range("K8").Value = 0.00867918216013561
$ws.Range("L8").Value = 0.008456762397614212
$ws.Range("M8").Value = 0.007720614698782578
$ws.Range("N8").Value = 0.00728699863393909
$ws.Range("O8").Value = 0.007019898960384205
$ws.Range("P8").Value = 0.006644134072997111
$ws.Range("Q8").Value = 0.006644134072997111
$ws.Range("R8").Value = 0.00663209500682782
$ws.Range("S8").Value = 0.006573959220548179
$ws.Range("T8").Value = 0.006462019485681618
$ws.Range("U8").Value = 0.006462019485681618
$ws.Range("V8").Value = 0.006413167032785945
$ws.Range("W8").Value = 0.006343545108804149
$ws.Range("X8").Value = 0.006343545108804149
$ws.Range("Y8").Value = 0.006294172610753865
$ws.Range("C9").Value = 0.6900010108947754
$ws.Range("E9").Value = 322.2840425619033
$ws.Range("F9").Value = 0.0150355307767265
$ws.Range("G9").Value = 0.01209376724620526
$ws.Range("H9").Value = 0.01061008389663061
$ws.Range("I9").Value = 0.009615049880918432
$ws.Range("J9").Value = 0.009105086099668702
$ws.Range("K9").Value = 0.008392642625032719
$ws.Range("L9").Value = 0.008099157051108996
$ws.Range("M9").Value = 0.007690605311864913
$ws.Range("N9").Value = 0.007155816308569913
$ws.Range("O9").Value = 0.006927488850893251
$ws.Range("P9").Value = 0.006838045140799084
$ws.Range("Q9").Value = 0.006714245834767967
$ws.Range("R9").Value = 0.006589411563768397
$ws.Range("S9").Value = 0.006554647185649728
$ws.Range("T9").Value = 0.006534087060830408
$ws.Range("U9").Value = 0.006469894904900344
$ws.Range("V9").Value = 0.006469894904900344
$ws.Range("W9").Value = 0.006282340010953279
$ws.Range("X9").Value = 0.006282340010953279
$ws.Range("Y9").Value = 0.006282340010953279
$ws.Range("C10").Value = 0.8150005340576172
$ws.Range("E10").Value = 333.1910624220709
$ws.Range("F10").Value = 0.01556858398318656
$ws.Range("G10").Value = 0.01110503744729848
$ws.Range("H10").Value = 0.01042244034827254
$ws.Range("I10").Value = 0.009121335591016988
$ws.Range("J10").Value = 0.008273834883937303
$ws.Range("K10").Value = 0.007361634589013379
$ws.Range("L10").Value = 0.007361634589013379
$ws.Range("M10").Value = 0.007361634589013379
$ws.Range("N10").Value = 0.007361634589013379
$ws.Range("O10").Value = 0.007222363960406712
$ws.Range("P10").Value = 0.00691511530516305
$ws.Range("Q10").Value = 0.006643131443058954
$ws.Range("R10").Value = 0.006643131443058954
$ws.Range("S10").Value = 0.006643131443058954
$ws.Range("T10").Value = 0.006643131443058954
$ws.Range("U10").Value = 0.006616975062492196
$ws.Range("V10").Value = 0.00651935508630903
$ws.Range("W10").Value = 0.00651935508630903
$ws.Range("X10").Value = 0.006501923750704688
$ws.Range("Y10").Value = 0.00649495248386103
$ws.Range("C11").Value = 0.890000581741333
$ws.Range("E11").Value = 331.8368499559328
$ws.Range("F11").Value = 0.01548982037978572
$ws.Range("G11").Value = 0.01178956995637096
$ws.Range("H11").Value = 0.01047149030301134
$ws.Range("I11").Value = 0.009481611048102183
$ws.Range("J11").Value = 0.00850349700955174
$ws.Range("K11").Value = 0.007833113933871894
$ws.Range("L11").Value = 0.007543143037210861
$ws.Range("M11").Value = 0.007356501877353227
$ws.Range("N11").Value = 0.006790472659713893
$ws.Range("O11").Value = 0.006790472659713893
$ws.Range("P11").Value = 0.006790472659713893
$ws.Range("Q11").Value = 0.00675995156066755
$ws.Range("R11").Value = 0.00675995156066755
$ws.Range("S11").Value = 0.00675995156066755
$ws.Range("T11").Value = 0.006596186250338887
$ws.Range("U11").Value = 0.006596186250338887
$ws.Range("V11").Value = 0.006563827516684257
$ws.Range("W11").Value = 0.006563827516684257
$ws.Range("X11").Value = 0.006563827516684257
$ws.Range("Y11").Value = 0.006503777979685526
